$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3479558.36
$ws.Range("C9").Value = 542146.99
$ws.Range("D9").Value = 4021705.35
$ws.Range("E9").Value = 13.4805248723654
$ws.Range("F9").Value = 86.51947512763461
$ws.Range("G9").Value = -47.60419668601735
$ws.Range("H9").Value = -37.16398560180627
$ws.Range("I9").Value = 34890
$ws.Range("J9").Value = 1474
$ws.Range("K9").Value = 36364
$ws.Range("L9").Value = 25104
$ws.Range("M9").Value = 160.2017746175908
$ws.Range("N9").Value = 9.373014477825636
